$wb = $excel.ActiveWorkbook

# Update "展览" sheet (rId1 / sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 152
$ws1.Range("F3").Value = 104

# Update "全部类型" sheet (rId4 / sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 152
$ws4.Range("F3").Value = 104
